$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended after the last existing row (row 25 -> row 26).
# Force column A to be stored as literal text (not auto-converted to a
# date serial number) by temporarily switching the cell to a text number
# format, then resetting the style so no extra formatting is left on it.
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = "12/20/2025"
$ws.Range("A26").Style = "Normal"

$ws.Range("B26").Value = 12040.8
$ws.Range("C26").Value = 0.206572700787405
$ws.Range("D26").Value = 0.793427299212595
$ws.Range("E26").Value = -138.84
$ws.Range("F26").Value = -28.63
$ws.Range("G26").Value = -20988.26
$ws.Range("H26").Value = -68.72
$ws.Range("I26").Value = -465.55
$ws.Range("J26").Value = -15.77
